# "ML model retrained with all data"
# Columns J (col 10) and K (col 11), rows 1-51, get new constant prediction
# weights from the retrained model: J = 0.5, K = 0.3 (previously J1/K1 held
# the shared-string labels "r"/"s" and J2:J51/K2:K51 held 1 / 0.5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J1:J51").Value = 0.5
$ws.Range("K1:K51").Value = 0.3

# Update the saved selection/view to match the new focus area (K column).
$null = $ws.Range("K1:K51").Select()
